$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.14676527569973
$ws.Range("C2").Value = 7.396519030095789
$ws.Range("D2").Value = 5.825949858745518
$ws.Range("F2").Value = 29.02658571147445
$ws.Range("G2").Value = 37.03474423979916
$ws.Range("H2").Value = 16.30466422911785
$ws.Range("I2").Value = 25.56292635420807
$ws.Range("K2").Value = 10.29699371270513
$ws.Range("L2").Value = 10.7503203325334
$ws.Range("N2").Value = 19.7779739676484
$ws.Range("B3").Value = 12.89401322144127
$ws.Range("C3").Value = 7.34115751824093
$ws.Range("D3").Value = 5.800802499780209
$ws.Range("F3").Value = 29.01737789483753
$ws.Range("G3").Value = 37.00580167057596
$ws.Range("H3").Value = 16.34274514529553
$ws.Range("I3").Value = 25.62641996580857
$ws.Range("K3").Value = 10.12230679072458
$ws.Range("L3").Value = 10.72461861042657
$ws.Range("N3").Value = 19.84055737437035
$ws.Range("B4").Value = 12.73930919593043
$ws.Range("C4").Value = 7.306336008312418
$ws.Range("D4").Value = 5.785002058158153
$ws.Range("F4").Value = 29.01906448875049
$ws.Range("G4").Value = 36.99879052136627
$ws.Range("H4").Value = 16.36885005515901
$ws.Range("I4").Value = 25.67029134008377
$ws.Range("K4").Value = 10.01574874331737
$ws.Range("L4").Value = 10.71100439704811
$ws.Range("N4").Value = 19.88078248782664
$ws.Range("B5").Value = 12.67647966989863
$ws.Range("C5").Value = 7.291941566876563
$ws.Range("D5").Value = 5.778474110809755
$ws.Range("F5").Value = 29.02159819719889
$ws.Range("G5").Value = 36.99863967698509
$ws.Range("H5").Value = 16.38017210570258
$ws.Range("I5").Value = 25.68939578504694
$ws.Range("K5").Value = 9.972560302759465
$ws.Range("L5").Value = 10.70600537531507
$ws.Range("N5").Value = 19.89762808397705
$ws.Range("B6").Value = 12.66606246722629
$ws.Range("C6").Value = 7.289539157815458
$ws.Range("D6").Value = 5.777384805889787
$ws.Range("F6").Value = 29.02213043032839
$ws.Range("G6").Value = 36.99877805343086
$ws.Range("H6").Value = 16.38209342072504
$ws.Range("I6").Value = 25.69264207650778
$ws.Range("K6").Value = 9.965404798425771
$ws.Range("L6").Value = 10.70520855059187
$ws.Range("N6").Value = 19.90045271165898
$ws.Range("B7").Value = 12.73846086857808
$ws.Range("C7").Value = 7.306142701916013
$ws.Range("D7").Value = 5.784914379121093
$ws.Range("F7").Value = 29.01909118330583
$ws.Range("G7").Value = 36.99877753044288
$ws.Range("H7").Value = 16.36899997936396
$ws.Range("I7").Value = 25.67054402618308
$ws.Range("K7").Value = 10.01516526035106
$ws.Range("L7").Value = 10.71093475113881
$ws.Range("N7").Value = 19.88100783543063
$ws.Range("B8").Value = 13.05956962300725
$ws.Range("C8").Value = 7.377604974568344
$ws.Range("D8").Value = 5.817354544365632
$ws.Range("F8").Value = 29.02188814685852
$ws.Range("G8").Value = 37.02253274224502
$ws.Range("H8").Value = 16.31722897862855
$ws.Range("I8").Value = 25.5838036107255
$ws.Range("K8").Value = 10.23665139060119
$ws.Range("L8").Value = 10.74101123931089
$ws.Range("N8").Value = 19.79918017342904
$ws.Range("B9").Value = 13.68920371570151
$ws.Range("C9").Value = 7.510984570172245
$ws.Range("D9").Value = 5.878056588540887
$ws.Range("F9").Value = 29.08553383864254
$ws.Range("G9").Value = 37.15435991836791
$ws.Range("H9").Value = 16.23734038386614
$ws.Range("I9").Value = 25.45256557100487
$ws.Range("K9").Value = 10.67399093764638
$ws.Range("L9").Value = 10.81698448479292
$ws.Range("N9").Value = 19.65292963791421
$ws.Range("B10").Value = 14.14668318951728
$ws.Range("C10").Value = 7.604628225906453
$ws.Range("D10").Value = 5.920803658189908
$ws.Range("F10").Value = 29.16756566078501
$ws.Range("G10").Value = 37.30285557633723
$ws.Range("H10").Value = 16.19187104050285
$ws.Range("I10").Value = 25.37995313214515
$ws.Range("K10").Value = 10.99379242631152
$ws.Range("L10").Value = 10.88286447254209
$ws.Range("N10").Value = 19.55406108211008
$ws.Range("B11").Value = 14.35268232055382
$ws.Range("C11").Value = 7.646236808759886
$ws.Range("D11").Value = 5.939832839892839
$ws.Range("F11").Value = 29.21247609250133
$ws.Range("G11").Value = 37.38150382367404
$ws.Range("H11").Value = 16.17406408142863
$ws.Range("I11").Value = 25.35211190550207
$ws.Range("K11").Value = 11.13827272452698
$ws.Range("L11").Value = 10.91494668328061
$ws.Range("N11").Value = 19.51092894345495
$ws.Range("B12").Value = 14.43030651757416
$ws.Range("C12").Value = 7.661846460236714
$ws.Range("D12").Value = 5.946977411873245
$ws.Range("F12").Value = 29.23056656527195
$ws.Range("G12").Value = 37.41286744490465
$ws.Range("H12").Value = 16.1677352247596
$ws.Range("I12").Value = 25.34231709571427
$ws.Range("K12").Value = 11.19278631251635
$ws.Range("L12").Value = 10.92739230630067
$ws.Range("N12").Value = 19.49485974616274
$ws.Range("B13").Value = 14.41360702283685
$ws.Range("C13").Value = 7.65849122929396
$ws.Range("D13").Value = 5.945441452690001
$ws.Range("F13").Value = 29.22662239627761
$ws.Range("G13").Value = 37.40604265809559
$ws.Range("H13").Value = 16.16907982605618
$ws.Range("I13").Value = 25.34439328858275
$ws.Range("K13").Value = 11.18105547904669
$ws.Range("L13").Value = 10.9246988362364
$ws.Range("N13").Value = 19.49830881103134
$ws.Range("B14").Value = 14.35907663071267
$ws.Range("C14").Value = 7.647523986352088
$ws.Range("D14").Value = 5.940421866258727
$ws.Range("F14").Value = 29.21394273545523
$ws.Range("G14").Value = 37.38405252881587
$ws.Range("H14").Value = 16.17353509772489
$ws.Range("I14").Value = 25.35129107951319
$ws.Range("K14").Value = 11.1427618600762
$ws.Range("L14").Value = 10.91596468617785
$ws.Range("N14").Value = 19.50960163612702
$ws.Range("B15").Value = 14.32562299657716
$ws.Range("C15").Value = 7.640787016389454
$ws.Range("D15").Value = 5.93733919013391
$ws.Range("F15").Value = 29.20631696516341
$ws.Range("G15").Value = 37.37078839256547
$ws.Range("H15").Value = 16.1763180418384
$ws.Range("I15").Value = 25.35561364122482
$ws.Range("K15").Value = 11.11927856181929
$ws.Range("L15").Value = 10.91065320013558
$ws.Range("N15").Value = 19.51655316644256
$ws.Range("B16").Value = 14.13317127806563
$ws.Range("C16").Value = 7.601888717287651
$ws.Range("D16").Value = 5.919551547760232
$ws.Range("F16").Value = 29.1647828018966
$ws.Range("G16").Value = 37.29793784312824
$ws.Range("H16").Value = 16.19309270998986
$ws.Range("I16").Value = 25.38187721520204
$ws.Range("K16").Value = 10.9843254986822
$ws.Range("L16").Value = 10.88080975055193
$ws.Range("N16").Value = 19.55691687293838
$ws.Range("B17").Value = 14.0145097191627
$ws.Range("C17").Value = 7.577769196517887
$ws.Range("D17").Value = 5.908531650791328
$ws.Range("F17").Value = 29.1412425391125
$ws.Range("G17").Value = 37.25607902904803
$ws.Range("H17").Value = 16.20412074589939
$ws.Range("I17").Value = 25.399319730948
$ws.Range("K17").Value = 10.90124089455918
$ws.Range("L17").Value = 10.86303807797659
$ws.Range("N17").Value = 19.58215011285565
$ws.Range("B18").Value = 13.9460645690482
$ws.Range("C18").Value = 7.563803312762356
$ws.Range("D18").Value = 5.902154164129223
$ws.Range("F18").Value = 29.12841813249021
$ws.Range("G18").Value = 37.23304845415499
$ws.Range("H18").Value = 16.21073457740152
$ws.Range("I18").Value = 25.40984067414247
$ws.Range("K18").Value = 10.8533618994404
$ws.Range("L18").Value = 10.85301562616146
$ws.Range("N18").Value = 19.59683720482418
$ws.Range("B19").Value = 13.9228593703936
$ws.Range("C19").Value = 7.559058865457214
$ws.Range("D19").Value = 5.899988179854004
$ws.Range("F19").Value = 29.12419909231547
$ws.Range("G19").Value = 37.22543066544464
$ws.Range("H19").Value = 16.21302040021769
$ws.Range("I19").Value = 25.41348672613048
$ws.Range("K19").Value = 10.83713697163815
$ws.Range("L19").Value = 10.84965664203136
$ws.Range("N19").Value = 19.60183985175004
$ws.Range("B20").Value = 14.02716213082251
$ws.Range("C20").Value = 7.58034641374669
$ws.Range("D20").Value = 5.909708798316113
$ws.Range("F20").Value = 29.14367446072509
$ws.Range("G20").Value = 37.26042686222875
$ws.Range("H20").Value = 16.20291876228526
$ws.Range("I20").Value = 25.39741237794005
$ws.Range("K20").Value = 10.91009520436587
$ws.Range("L20").Value = 10.86490932132816
$ws.Range("N20").Value = 19.57944603252158
$ws.Range("B21").Value = 14.37510454197813
$ws.Range("C21").Value = 7.650749343637451
$ws.Range("D21").Value = 5.941897917421235
$ws.Range("F21").Value = 29.21763771458297
$ws.Range("G21").Value = 37.39046877583941
$ws.Range("H21").Value = 16.17221523009872
$ws.Range("I21").Value = 25.349244716655
$ws.Range("K21").Value = 11.15401541756667
$ws.Range("L21").Value = 10.91852212157302
$ws.Range("N21").Value = 19.5062775033413
$ws.Range("B22").Value = 14.60023761679149
$ws.Range("C22").Value = 7.695905248543456
$ws.Range("D22").Value = 5.962577014626622
$ws.Range("F22").Value = 29.27228988311241
$ws.Range("G22").Value = 37.48466693627351
$ws.Range("H22").Value = 16.15456345750894
$ws.Range("I22").Value = 25.32212526448319
$ws.Range("K22").Value = 11.31225406710722
$ws.Range("L22").Value = 10.95528750328176
$ws.Range("N22").Value = 19.45999597198232
$ws.Range("B23").Value = 14.48031250348581
$ws.Range("C23").Value = 7.671884426658342
$ws.Range("D23").Value = 5.951573456350989
$ws.Range("F23").Value = 29.24254638311449
$ws.Range("G23").Value = 37.43355445236151
$ws.Range("H23").Value = 16.16376344265793
$ws.Range("I23").Value = 25.33619990152444
$ws.Range("K23").Value = 11.22792414380793
$ws.Range("L23").Value = 10.93550961525865
$ws.Range("N23").Value = 19.4845569019814
$ws.Range("B24").Value = 14.02144267033308
$ws.Range("C24").Value = 7.579181562683801
$ws.Range("D24").Value = 5.909176740623708
$ws.Range("F24").Value = 29.1425727792095
$ws.Range("G24").Value = 37.25845798297581
$ws.Range("H24").Value = 16.20346132663446
$ws.Range("I24").Value = 25.39827315663663
$ws.Range("K24").Value = 10.90609251729648
$ws.Range("L24").Value = 10.864062725153
$ws.Range("N24").Value = 19.58066798592012
$ws.Range("B25").Value = 13.51941344539757
$ws.Range("C25").Value = 7.475649330567028
$ws.Range("D25").Value = 5.86195468487805
$ws.Range("F25").Value = 29.06210330794293
$ws.Range("G25").Value = 37.10959506938585
$ws.Range("H25").Value = 16.2566323109937
$ws.Range("I25").Value = 25.48389591799637
$ws.Range("K25").Value = 10.5557080769031
$ws.Range("L25").Value = 10.79464118277395
$ws.Range("N25").Value = 19.69098128778839
